$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/date/measurement cells: direct .Value assignment is safe
# (Excel COM does not coerce these into numbers because of the unit
# suffix, e.g. " hPa", " mm", " °C", or the date-time format).
$ws.Range("E2").Value = "2026-02-08 22:18:42"
$ws.Range("E3").Value = "2026-02-08 22:18:45"
$ws.Range("E4").Value = "2026-02-08 22:18:48"
$ws.Range("J4").Value = "1002.7 hPa"
$ws.Range("O4").Value = "10.1 °C"
$ws.Range("E5").Value = "2026-02-08 22:18:50"
$ws.Range("E6").Value = "2026-02-08 22:18:53"
$ws.Range("J6").Value = "1002.6 hPa"
$ws.Range("E7").Value = "2026-02-08 22:18:56"
$ws.Range("J7").Value = "1002.9 hPa"
$ws.Range("E8").Value = "2026-02-08 22:18:59"
$ws.Range("J8").Value = "1002.8 hPa"
$ws.Range("O8").Value = "9.3 °C"
$ws.Range("E9").Value = "2026-02-08 22:19:02"
$ws.Range("E10").Value = "2026-02-08 22:19:04"
$ws.Range("I10").Value = "3.9 mm"
$ws.Range("E11").Value = "2026-02-08 22:19:07"
$ws.Range("E12").Value = "2026-02-08 22:19:10"
$ws.Range("E13").Value = "2026-02-08 22:19:12"
$ws.Range("J13").Value = "1004.1 hPa"
$ws.Range("E14").Value = "2026-02-08 22:19:15"
$ws.Range("E15").Value = "2026-02-08 22:19:18"
$ws.Range("E16").Value = "2026-02-08 22:19:20"
$ws.Range("E17").Value = "2026-02-08 22:19:22"
$ws.Range("E18").Value = "2026-02-08 22:19:25"
$ws.Range("I18").Value = "0.7 mm"
$ws.Range("J18").Value = "1002.9 hPa"
$ws.Range("E19").Value = "2026-02-08 22:19:28"
$ws.Range("E20").Value = "2026-02-08 22:19:30"
$ws.Range("I20").Value = "10.0 mm"
$ws.Range("O20").Value = "-4.4 °C"
$ws.Range("E21").Value = "2026-02-08 22:19:33"
$ws.Range("J21").Value = "1003.6 hPa"
$ws.Range("O21").Value = "5.3 °C"
$ws.Range("E22").Value = "2026-02-08 22:19:35"
$ws.Range("E23").Value = "2026-02-08 22:19:38"
$ws.Range("E24").Value = "2026-02-08 22:19:41"
$ws.Range("J24").Value = "1004.1 hPa"
$ws.Range("E25").Value = "2026-02-08 22:19:44"
$ws.Range("E26").Value = "2026-02-08 22:19:46"
$ws.Range("J26").Value = "1002.0 hPa"
$ws.Range("E27").Value = "2026-02-08 22:19:49"
$ws.Range("E28").Value = "2026-02-08 22:19:52"
$ws.Range("J28").Value = "1002.5 hPa"
$ws.Range("E29").Value = "2026-02-08 22:19:55"
$ws.Range("I29").Value = "3.1 mm"
$ws.Range("O29").Value = "10.5 °C"
$ws.Range("E30").Value = "2026-02-08 22:19:57"
$ws.Range("J30").Value = "1002.9 hPa"
$ws.Range("E31").Value = "2026-02-08 22:20:00"
$ws.Range("J31").Value = "1002.1 hPa"
$ws.Range("N31").Value = "7.3 °C 21:55 TU"
$ws.Range("O31").Value = "9.5 °C"
$ws.Range("E32").Value = "2026-02-08 22:20:03"
$ws.Range("E33").Value = "2026-02-08 22:20:06"
$ws.Range("J33").Value = "1003.7 hPa"
$ws.Range("E34").Value = "2026-02-08 22:20:09"
$ws.Range("K34").Value = "12.5 MJ/m2"
$ws.Range("E35").Value = "2026-02-08 22:20:11"
$ws.Range("J35").Value = "1005.0 hPa"
$ws.Range("O35").Value = "4.0 °C"
$ws.Range("E36").Value = "2026-02-08 22:20:14"
$ws.Range("J36").Value = "1002.9 hPa"
$ws.Range("E37").Value = "2026-02-08 22:20:17"
$ws.Range("J37").Value = "1003.8 hPa"
$ws.Range("E38").Value = "2026-02-08 22:20:20"
$ws.Range("E39").Value = "2026-02-08 22:20:22"
$ws.Range("E40").Value = "2026-02-08 22:20:25"
$ws.Range("J40").Value = "1004.2 hPa"
$ws.Range("O40").Value = "5.6 °C"
$ws.Range("E41").Value = "2026-02-08 22:20:28"
$ws.Range("J41").Value = "1002.9 hPa"
$ws.Range("O41").Value = "12.2 °C"
$ws.Range("E42").Value = "2026-02-08 22:20:31"
$ws.Range("O42").Value = "10.6 °C"
$ws.Range("E43").Value = "2026-02-08 22:20:33"
$ws.Range("O43").Value = "7.2 °C"
$ws.Range("E44").Value = "2026-02-08 22:20:36"
$ws.Range("I44").Value = "2.4 mm"
$ws.Range("E45").Value = "2026-02-08 22:20:38"
$ws.Range("J45").Value = "1005.1 hPa"
$ws.Range("E46").Value = "2026-02-08 22:20:41"
$ws.Range("J46").Value = "1004.6 hPa"

# Percentage-looking cells (e.g. "71%") are stored as plain text in the
# source data, but a direct .Value assignment makes Excel COM interpret
# them as a numeric percentage (0.71) instead of keeping literal text.
# Work around this by writing the text into a scratch cell formatted as
# Text ("@"), copying it, and pasting-special (values only) onto the
# target cell - this preserves both the literal string and the original
# cell style of the destination.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "71%"
$scratch.Copy()
$ws.Range("H4").PasteSpecial(-4163)
$scratch.NumberFormat = "@"
$scratch.Value = "71%"
$scratch.Copy()
$ws.Range("H18").PasteSpecial(-4163)
$scratch.NumberFormat = "@"
$scratch.Value = "80%"
$scratch.Copy()
$ws.Range("H21").PasteSpecial(-4163)
$scratch.NumberFormat = "@"
$scratch.Value = "76%"
$scratch.Copy()
$ws.Range("H36").PasteSpecial(-4163)
$scratch.NumberFormat = "@"
$scratch.Value = "78%"
$scratch.Copy()
$ws.Range("H38").PasteSpecial(-4163)
$scratch.NumberFormat = "@"
$scratch.Value = "86%"
$scratch.Copy()
$ws.Range("H39").PasteSpecial(-4163)
$scratch.NumberFormat = "@"
$scratch.Value = "67%"
$scratch.Copy()
$ws.Range("H41").PasteSpecial(-4163)
$scratch.Clear()
